# Updated cryptos list on Fri Nov 24 11:35:13 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = "37.801.19"
$ws.Range("E2").Value = "  +1.22%  "
$ws.Range("D3").Value = "2.105.33"
$ws.Range("E3").Value = "  +2.41%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'235.18"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.87%  "
$ws.Range("D6").Value = "'0.624"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.16%  "
$ws.Range("D7").Value = "'58.39"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.00%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  +2.19%  "
$ws.Range("E10").Value = "  +2.34%  "
$ws.Range("E11").Value = "  +1.13%  "
$ws.Range("D12").Value = "2.417.83"
$ws.Range("E12").Value = "  +2.68%  "
$ws.Range("D13").Value = "'14.54"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.18%  "
$ws.Range("D14").Value = "'21.25"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.87%  "
$ws.Range("D15").Value = "'0.785"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.49%  "
$ws.Range("D16").Value = "'5.25"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.66%  "
$ws.Range("D17").Value = "2.131.78"
$ws.Range("E17").Value = "  +3.90%  "
$ws.Range("D18").Value = "37.702.39"
$ws.Range("E18").Value = "  +0.55%  "
$ws.Range("D19").Value = "'6.27"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.24%  "
$ws.Range("D20").Value = "'70.09"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.29%  "
$ws.Range("E21").Value = "  +1.33%  "
$ws.Range("D22").Value = "'227.80"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.82%  "
$ws.Range("E23").Value = "  -0.10%  "
$ws.Range("E24").Value = "  +0.94%  "
$ws.Range("E25").Value = "  +0.05%  "
$ws.Range("D26").Value = "'167.94"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.61%  "
$ws.Range("E27").Value = "  +0.88%  "
$ws.Range("E28").Value = "  +3.12%  "
$ws.Range("E29").Value = "  -3.78%  "
$ws.Range("D30").Value = "'19.48"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.86%  "
$ws.Range("E31").Value = "  +0.63%  "
$ws.Range("D32").Value = "'4.68"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.56%  "
$ws.Range("D33").Value = "'2.60"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.03%  "
$ws.Range("D34").Value = "'0.0623"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.39%  "
$ws.Range("E35").Value = "  +0.23%  "
$ws.Range("E36").Value = "  +5.56%  "
$ws.Range("D37").Value = "'1.77"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.78%  "
$ws.Range("E38").Value = "  +0.23%  "
$ws.Range("D39").Value = "'5.61"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.38%  "
$ws.Range("E40").Value = "  -0.14%  "
$ws.Range("E41").Value = "  +2.02%  "
$ws.Range("D42").Value = "'97.53"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.77%  "
$ws.Range("D43").Value = "1.474.46"
$ws.Range("E43").Value = "  +1.18%  "
$ws.Range("E44").Value = "  +1.13%  "
$ws.Range("E45").Value = "  -0.75%  "
$ws.Range("D46").Value = "'4.22"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -9.69%  "
$ws.Range("D49").Value = "'15.58"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.55%  "
$ws.Range("E50").Value = "  +3.12%  "
$ws.Range("D51").Value = "2.301.98"
$ws.Range("E51").Value = "  +2.63%  "

# Row 47/48: ARBITRUM and FraxShare swap positions with updated values
$ws.Range("B47").Value = "ARBITRUM"
$ws.Range("C47").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D47").Value = "'1.05"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.23%  "
$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").Value = "'7.47"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.27%  "
